# "Generate Report for Handoff"
#
# The localization report moves from "In Translation" to "Ready for
# handoff": the Status text + the associated "Latest Handoff Datetime" /
# "Latest HO Xliff Generate Date" timestamps are refreshed on all three
# sheets, and the Status column is widened to fit the new, longer text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status: "In Translation" -> "Ready for handoff" -----------------------
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value     = "Ready for handoff"
$dede.Range("C2").Value     = "Ready for handoff"

# --- Handoff timestamps, refreshed for the new handoff -----------------
$overview.Range("G2").Value = "2016-08-30 10:44:54"
$dede.Range("H2").Value     = "2016-08-30 10:44:54"
$zhcn.Range("H2").Value     = "2016-08-30 10:44:49"

# --- Widen the Status column(s) so the longer "Ready for handoff" text fits
$overview.Range("E1").ColumnWidth = 16.333333333333336
$overview.Range("F1").ColumnWidth = 16.333333333333336
$zhcn.Range("C1").ColumnWidth     = 16.333333333333336
$dede.Range("C1").ColumnWidth     = 16.333333333333336
